# update data collection model
#
# DataCollectionGroup: insert "implementationOption" / "scenario" columns
#   before "domain" (old H, new J), add a list validation on the new
#   "implementationOption" column (H2:H1048576).
#
# DataCollectionItem: insert "variableName" right after "name"; drop the old
#   "isNonStandard" / "dataCollectionInstrumentItem" columns; move "codelist"
#   to just after "displayHidden"; append a new "sdtmAnnotation" column at
#   the end; move the dataType list validation from column I to column H
#   (its new position).
#
# SDTMTarget: drop "sdtmAnnotation" (it now lives on DataCollectionItem) and
#   rename the remaining "sdtmTargetMapping" column from C to B.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# DataCollectionGroup
# ---------------------------------------------------------------------
$wsGroup = $wb.Worksheets.Item("DataCollectionGroup")

# Make room for the two new columns right before "domain" (H).
$wsGroup.Range("H1:I1").EntireColumn.Insert()

$wsGroup.Range("H1").Value = "implementationOption"
$wsGroup.Range("I1").Value = "scenario"

$implOptionRange = $wsGroup.Range("H2:H1048576")
$implOptionRange.Validation.Add(3, 1, 1, '"Horizontal,Vertical"')
$implOptionRange.Validation.ShowInput = $false
$implOptionRange.Validation.ShowError = $false

# ---------------------------------------------------------------------
# DataCollectionItem
# ---------------------------------------------------------------------
$wsItem = $wb.Worksheets.Item("DataCollectionItem")

$itemHeaders = @(
    "name",
    "variableName",
    "dataElementConceptId",
    "questionText",
    "prompt",
    "orderNumber",
    "mandatoryVariable",
    "dataType",
    "length",
    "significantDigits",
    "displayHidden",
    "codelist",
    "valueList",
    "listType",
    "prepopulatedValue",
    "sdtmTarget",
    "sdtmAnnotation"
)
for ($i = 0; $i -lt $itemHeaders.Length; $i++) {
    $wsItem.Cells.Item(1, $i + 1).Value = $itemHeaders[$i]
}

# "dataType" validation used to live on I (old layout); it now belongs on H.
$oldDataTypeRange = $wsItem.Range("I2:I1048576")
$listTypeRange = $wsItem.Range("N2:N1048576")

$oldDataTypeRange.Validation.Delete()
$listTypeRange.Validation.Delete()

$dataTypeRange = $wsItem.Range("H2:H1048576")
$dataTypeRange.Validation.Add(3, 1, 1, '"decimal,float,integer,text,date,time"')
$dataTypeRange.Validation.ShowInput = $false
$dataTypeRange.Validation.ShowError = $false

$listTypeRange.Validation.Add(3, 1, 1, '"Radio,Dropdown,DropdownMultiSelect,Checkbox,Text,Date,Time,DateTime"')
$listTypeRange.Validation.ShowInput = $false
$listTypeRange.Validation.ShowError = $false

# ---------------------------------------------------------------------
# SDTMTarget
# ---------------------------------------------------------------------
$wsSdtm = $wb.Worksheets.Item("SDTMTarget")

$wsSdtm.Range("B1").Value = "sdtmTargetMapping"
$wsSdtm.Range("C1").EntireColumn.Delete()
